{"js": "// The original paragraph reads:\n//   \"Para el servidor nos interesa utilizar el Framework Spring ...\"\n// and needs to become:\n//   \"Primero para el desarrollo del BackEnd utilizaremos Java, acompa\u00f1ado de\n//    lo siguiente para el correcto funcionamiento de este y para el servidor\n//    nos interesa utilizar el Framework Spring ...\"\n// i.e. a chunk of text is inserted right after the leading \"P\" (before\n// \"ara el servidor...\"). Find the unique substring \"ara el servidor\" and\n// insert the new text just before it.\n\nconst body = context.document.body;\n\nconst searchResults = body.search(\"ara el servidor\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Expected text \"ara el servidor\" not found in document body.');\n}\n\nconst target = searchResults.items[0];\nconst insertion =\n  \"rimero para el desarrollo del BackEnd utilizaremos Java, acompa\u00f1ado de lo siguiente para el correcto funcionamiento de este y p\";\n\ntarget.insertText(insertion, \"Before\");\nawait context.sync();\n", "ps1": "# The paragraph under \"Tecnolog\u00edas Empleadas\" originally starts:\n#   \"Para el servidor nos interesa utilizar el Framework Spring ...\"\n# It needs to become:\n#   \"Primero para el desarrollo del BackEnd utilizaremos Java, acompa\u00f1ado de\n#    lo siguiente para el correcto funcionamiento de este y para el servidor\n#    nos interesa utilizar el Framework Spring ...\"\n# Everything from \"servidor nos interesa ...\" to the end of the paragraph is\n# unchanged, so a single Find & Replace on the unique lead-in phrase\n# \"Para el servidor\" reproduces the edit exactly.\n\n$d = $word.ActiveDocument\n\n$findText = \"Para el servidor\"\n$replaceText = \"Primero para el desarrollo del BackEnd utilizaremos Java, acompa\u00f1ado de lo siguiente para el correcto funcionamiento de este y para el servidor\"\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n"}
